$d = $word.ActiveDocument

$p0 = @'
<w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="44"/></w:rPr><w:t>Medicine in the Realm of Human Health: Unmasking the Science of Life</w:t></w:r></w:p>
'@
$p1 = @'
<w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:t>Dr</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:t xml:space="preserve"> Louisa Harrison</w:t></w:r></w:p>
'@
$p2 = @'
<w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>louisa</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>harrison@hghs</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>edu</w:t></w:r></w:p>
'@
$p3 = @'
<w:p/>
'@
$p4 = @'
<w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>The vast discipline of medicine is like an intricate symphony, where diverse elements collaboratively compose a melody of health and life</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> It requires artistic sensitivity, robust understanding of science, and relentless pursuit of empathy and compassion</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Medicine not only enriches our knowledge of the human body and its intricacies but also compels us to explore the depths of our own understanding beyond organic processes and delve into the profound impact it has on our lives and relationships</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> It invites us to unveil the extraordinary workings of human existence and to question the nature of life itself</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/><w:t>The acquisition of medical knowledge is a continuous process that challenges us to understand intricate systems and mechanisms beyond textbooks and theory</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Harnessing the expertise of scientists, researchers, and practitioners, we immerse ourselves in the beauty of integrating anatomy, physiology, biochemistry, and genetics to create a comprehensive map of the human body</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> This unparalleled knowledge opens doors to diagnosis, treatment, and prevention of diseases, empowering us to create avenues of health and vitality for individuals and communities</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/><w:t>With its inherent nature of service and empathy, medicine extends beyond laboratories and wards, weaving its way into the very fabric of human experience</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> It teaches us the art of listening, transcending words to understand the language of body and spirit</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> This understanding enables us to navigate illness not simply as anomalies, but as opportunities for nurturing and care, guiding patients toward a state of holistic well-being that encompasses physical, mental, and spiritual health</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r></w:p>
'@
$p5 = @'
<w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="28"/></w:rPr><w:t>Summary</w:t></w:r></w:p>
'@
$p6 = @'
<w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t>In this essay, we embarked on a journey into the realm of medicine, an art and science intertwined</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> Medicine's prowess extends from understanding our bodies' complexity to comprehending our existence's essence</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> It invites us to contemplate the miracle of life and inspires us to serve as healers, listeners, and guardians of health</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> Above all, it reminds us that medicine is not merely a realm of treatment, but a harmonious union of science, compassion, and human connectivity that transforms lives and its communities</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r></w:p>
'@
$p7 = @'
<w:p/>
'@

$paras = @($p0, $p1, $p2, $p3, $p4, $p5, $p6)

for ($i = 0; $i -le 6; $i++) {
    $xml = $paras[$i]
    $para = $d.Paragraphs.Item($i + 1)
    $para.Range.InsertXML($xml)
}

# Append a new empty trailing paragraph (matches diff's added <w:p/> before sectPr)
$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.InsertXML($p7)

Write-Host "Paragraph count:" $d.Paragraphs.Count
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Host $i ":" $d.Paragraphs.Item($i).Range.Text
}
